# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (fund-level holdings detail) right
#    before the "总计" (totals) sheet.
# 2. Insert a new top row into "总计" summarizing the 2022-Q1 quarter and
#    renumber the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" sheet with per-fund holdings
# ---------------------------------------------------------------------

# NOTE: worksheet references are positional anchors here -- once
# Worksheets.Add()/rename reshuffles sheet order, a handle obtained
# *before* the structural change can silently resolve to a different
# sheet afterwards. Only use $totalSheetBefore to position the insert,
# then re-fetch a fresh "总计" handle by name for Part 2 below.
$totalSheetBefore = $wb.Worksheets.Item("总计")
$template          = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Copy header (row1 B:H) and index-column (A2:A17) formatting from the
# previous quarter sheet so the new sheet matches the house style exactly.
$template.Range("B1:H1").Copy() | Out-Null
$newSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$template.Range("A2:A13").Copy() | Out-Null
$newSheet.Range("A2:A17").PasteSpecial(-4122) | Out-Null

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns D/E/F/G hold percentage-looking / decimal-looking numbers that
# are actually stored as plain text in the source data (e.g. "34.40",
# "011338" for codes) -- force text format *before* assigning values so
# leading/trailing zeros survive instead of being normalised as numbers.
$newSheet.Range("B2:G17").NumberFormat = "@"

$funds = @(
    @("011338", "兴全合远两年持有期混合A",                           "34.40", "88.75", "3.03", "1.0423", 9),
    @("009896", "广发港股通成长精选股票A",                           "27.73", "89.63", "3.49", "0.9678", 10),
    @("513060", "博时恒生医疗保健ETF（QDII）",                       "23.35", "98.89", "4.09", "0.9550", 6),
    @("008269", "大成睿享混合A",                                     "17.69", "65.25", "3.93", "0.6952", 4),
    @("009897", "广发港股通成长精选股票C",                           "6.49",  "89.63", "3.49", "0.2265", 10),
    @("011834", "大成投资严选六个月持有期混合型证券投资基金A",       "3.88",  "84.63", "5.28", "0.2049", 7),
    @("013463", "大成致远优势一年持有期混合A",                       "4.01",  "60.15", "4.87", "0.1953", 4),
    @("090013", "大成竞争优势混合",                                  "3.87",  "72.16", "4.58", "0.1772", 3),
    @("513700", "鹏华中证港股通医药卫生综合交易型开放式指数证券投资基金", "3.24", "93.11", "3.48", "0.1128", 7),
    @("008270", "大成睿享混合C",                                     "2.87",  "65.25", "3.93", "0.1128", 4),
    @("011651", "招商港股通核心精选股票A",                           "2.81",  "81.27", "2.93", "0.0823", 6),
    @("159892", "华夏恒生香港上市生物科技ETF（QDII）",               "1.51",  "99.03", "4.47", "0.0675", 6),
    @("011339", "兴全合远两年持有期混合C",                           "1.53",  "88.75", "3.03", "0.0464", 9),
    @("011652", "招商港股通核心精选股票C",                           "0.94",  "81.27", "2.93", "0.0275", 6),
    @("011835", "大成投资严选六个月持有期混合型证券投资基金C",       "0.30",  "84.63", "5.28", "0.0158", 7),
    @("013464", "大成致远优势一年持有期混合C",                       "0.17",  "60.15", "4.87", "0.0083", 4)
)

for ($i = 0; $i -lt $funds.Count; $i++) {
    $r    = $i + 2
    $fund = $funds[$i]

    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $fund[0]
    $newSheet.Cells.Item($r, 3).Value = $fund[1]
    $newSheet.Cells.Item($r, 4).Value = $fund[2]
    $newSheet.Cells.Item($r, 5).Value = $fund[3]
    $newSheet.Cells.Item($r, 6).Value = $fund[4]
    $newSheet.Cells.Item($r, 7).Value = $fund[5]
    $newSheet.Cells.Item($r, 8).Value = $fund[6]
}

# ---------------------------------------------------------------------
# Part 2: prepend the 2022-Q1 summary row on the "总计" sheet
# ---------------------------------------------------------------------

# Re-fetch by name: the sheet collection was restructured above, so any
# handle grabbed beforehand is no longer trustworthy.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()

# Re-use the existing index-column formatting for the newly inserted cell.
$totalSheet.Range("A3").Copy() | Out-Null
$totalSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 4.94

# Renumber the 0-based index column now that a row was prepended.
$lastRow = $totalSheet.Cells.Item($totalSheet.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
